# Applies the "Ran a few tests" edit to PriorsTable2.xlsx:
#  - Switch workbook calculation to manual (calcPr calcMode="manual")
#  - Bump the E2:E22 prior-std-dev column from 0.05 to 0.15
#  - Update the sheet's selection to E2:E22 and drop the scrolled-down
#    topLeftCell that was previously in view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Calculation mode -----------------------------------------------------
# xlCalculationManual = -4135
$excel.Calculation = -4135

# --- Cell edits -------------------------------------------------------------
# E2:E22 all move from 0.05 to 0.15
$ws.Range("E2:E22").Value = 0.15

# --- View / selection ---------------------------------------------------
# Select E2:E22 so the saved sheetView reflects the new selection instead
# of the stale F46 one (the runtime always anchors the active cell at the
# top-left of whatever is selected, so this lands on E2 rather than E3,
# but E2:E22 is the best achievable match for the recorded sqref).
$ws.Range("E2:E22").Select()
